# Natmi following Dr Hou advice
#
# The LR-pair result table (Qrfp -> P2ry14) is recomputed for three
# sending/target cell-type clusters (ECs, FAPs, sCs) instead of the
# original two (FAPs, sCs), producing the full 3x3 cross of
# sending-cluster x target-cluster rows (9 data rows instead of 2),
# each with refreshed NATMI metric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=Sending cluster, B=Ligand symbol, C=Receptor symbol,
# D=Target cluster, then E..T are the 16 numeric NATMI metric columns.
$rows = @(
    ,@("ECs","Qrfp","P2ry14","ECs",2,0.6666666666666666,0.3347443333333333,1.004233,0.162471248533204,0.162471248533204,2,0.6666666666666666,0.8985673333333333,2.695702,0.1174926698436198,0.1174926698436198,0.3007903229517777,2.707112906566,0.01908918076299244,0.01908918076299244)
    ,@("ECs","Qrfp","P2ry14","FAPs",2,0.6666666666666666,0.3347443333333333,1.004233,0.162471248533204,0.162471248533204,3,1,4.590669333333333,13.772008,0.6002555137873887,0.6002555137873886,1.536700545540444,13.830304909864,0.0975242627639769,0.09752426276397687)
    ,@("ECs","Qrfp","P2ry14","sCs",2,0.6666666666666666,0.3347443333333333,1.004233,0.162471248533204,0.162471248533204,3,1,2.158622,6.475866,0.2822518163689915,0.2822518163689915,0.7225864823086665,6.503278340777999,0.04585780500623467,0.04585780500623467)
    ,@("FAPs","Qrfp","P2ry14","ECs",3,1,0.9771596666666666,2.931479,0.4742734536495697,0.4742734536495697,2,0.6666666666666666,0.8985673333333333,2.695702,0.1174926698436198,0.1174926698436198,0.8780437559175555,7.902393803258,0.05572365430524223,0.05572365430524222)
    ,@("FAPs","Qrfp","P2ry14","FAPs",3,1,0.9771596666666666,2.931479,0.4742734536495697,0.4742734536495697,3,1,4.590669333333333,13.772008,0.6002555137873887,0.6002555137873886,4.485816915536889,40.372352239832,0.2846852555961418,0.2846852555961417)
    ,@("FAPs","Qrfp","P2ry14","sCs",3,1,0.9771596666666666,2.931479,0.4742734536495697,0.4742734536495697,3,1,2.158622,6.475866,0.2822518163689915,0.2822518163689915,2.109318353979333,18.983865185814,0.1338645437481857,0.1338645437481857)
    ,@("sCs","Qrfp","P2ry14","ECs",2,0.6666666666666666,0.7484256666666665,2.245277,0.3632552978172263,0.3632552978172263,2,0.6666666666666666,0.8985673333333333,2.695702,0.1174926698436198,0.1174926698436198,0.6725108554948888,6.052597699453999,0.04267983477538517,0.04267983477538517)
    ,@("sCs","Qrfp","P2ry14","FAPs",2,0.6666666666666666,0.7484256666666665,2.245277,0.3632552978172263,0.3632552978172263,3,1,4.590669333333333,13.772008,0.6002555137873887,0.6002555137873886,3.435774756246222,30.921972806216,0.21804599542727,0.21804599542727)
    ,@("sCs","Qrfp","P2ry14","sCs",2,0.6666666666666666,0.7484256666666665,2.245277,0.3632552978172263,0.3632552978172263,3,1,2.158622,6.475866,0.2822518163689915,0.2822518163689915,1.615568109431333,14.540112984882,0.1025294676145711,0.1025294676145711)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    for ($i = 0; $i -lt 16; $i++) {
        $ws.Cells.Item($r, 5 + $i).Value = $row[4 + $i]
    }
    $r++
}
